# Refresh computed leve-profit columns (H:N) across the affected ALC/ARM/CRP/CUL/GSM/LTW/WVR rows
# per the scheduled-runner market-price pull.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 98
$ws.Range("H98").Value = 1101.3158
$ws.Range("I98").Value = 1133.6666
$ws.Range("J98").Value = 519
$ws.Range("K98").Value = 1133.6666
$ws.Range("L98").Value = 519
$ws.Range("M98").Value = 364.3334
$ws.Range("N98").Value = -3515

# Row 116
$ws.Range("H116").Value = 6366.6665
$ws.Range("I116").Value = 3550
$ws.Range("K116").Value = 3550
$ws.Range("M116").Value = -108

# Row 122
$ws.Range("H122").Value = 1101.3158
$ws.Range("I122").Value = 1133.6666
$ws.Range("J122").Value = 519
$ws.Range("K122").Value = 3400.9998
$ws.Range("L122").Value = 1557
$ws.Range("M122").Value = -950.9998
$ws.Range("N122").Value = -6457

# Row 138
$ws.Range("H138").Value = 17243876
$ws.Range("I138").Value = 62501456
$ws.Range("J138").Value = 2893.8572
$ws.Range("K138").Value = 187504368
$ws.Range("L138").Value = 8681.5716
$ws.Range("M138").Value = -187499228
$ws.Range("N138").Value = -18961.5716

$ws = $wb.Worksheets.Item("ARM")
# Row 32
$ws.Range("H32").Value = 5251.094
$ws.Range("I32").Value = 5251.094
$ws.Range("K32").Value = 5251.094
$ws.Range("M32").Value = -4964.094

# Row 74
$ws.Range("H74").Value = 2464.4285
$ws.Range("I74").Value = 1869.826
$ws.Range("K74").Value = 1869.826
$ws.Range("M74").Value = -995.826

# Row 77
$ws.Range("H77").Value = 2464.4285
$ws.Range("I77").Value = 1869.826
$ws.Range("K77").Value = 9349.130000000001
$ws.Range("M77").Value = -4981.130000000001

$ws = $wb.Worksheets.Item("CRP")
# Row 31
$ws.Range("H31").Value = 2106.4211
$ws.Range("I31").Value = 2106.4211
$ws.Range("J31").Value = 0
$ws.Range("K31").Value = 2106.4211
$ws.Range("L31").Value = 0
$ws.Range("M31").Value = -1811.4211
$ws.Range("N31").ClearContents()

# Row 34
$ws.Range("H34").Value = 2106.4211
$ws.Range("I34").Value = 2106.4211
$ws.Range("J34").Value = 0
$ws.Range("K34").Value = 2106.4211
$ws.Range("L34").Value = 0
$ws.Range("M34").Value = -1904.4211
$ws.Range("N34").ClearContents()

# Row 132
$ws.Range("H132").Value = 32419.666
$ws.Range("I132").Value = 17500
$ws.Range("J132").Value = 36682.43
$ws.Range("K132").Value = 52500
$ws.Range("L132").Value = 110047.29
$ws.Range("M132").Value = -49970
$ws.Range("N132").Value = -115107.29

$ws = $wb.Worksheets.Item("CUL")
# Row 26
$ws.Range("H26").Value = 416.46667
$ws.Range("J26").Value = 855.1667
$ws.Range("L26").Value = 2565.5001
$ws.Range("N26").Value = -3141.5001

# Row 50
$ws.Range("H50").Value = 2662.2666
$ws.Range("I50").Value = 489.6
$ws.Range("J50").Value = 3748.6
$ws.Range("K50").Value = 1468.8
$ws.Range("L50").Value = 11245.8
$ws.Range("M50").Value = -987.8000000000002
$ws.Range("N50").Value = -12207.8

# Row 53
$ws.Range("H53").Value = 2662.2666
$ws.Range("I53").Value = 489.6
$ws.Range("J53").Value = 3748.6
$ws.Range("K53").Value = 1468.8
$ws.Range("L53").Value = 11245.8
$ws.Range("M53").Value = -987.8000000000002
$ws.Range("N53").Value = -12207.8

# Row 68
$ws.Range("H68").Value = 1369.5714
$ws.Range("J68").Value = 900
$ws.Range("L68").Value = 2700
$ws.Range("N68").Value = -4322

# Row 71
$ws.Range("H71").Value = 1369.5714
$ws.Range("J71").Value = 900
$ws.Range("L71").Value = 8100
$ws.Range("N71").Value = -16212

# Row 82
$ws.Range("H82").Value = 12151.875
$ws.Range("J82").Value = 12304.667
$ws.Range("L82").Value = 36914.001
$ws.Range("N82").Value = -37726.001

# Row 85
$ws.Range("H85").Value = 12151.875
$ws.Range("J85").Value = 12304.667
$ws.Range("L85").Value = 36914.001
$ws.Range("N85").Value = -39722.001

# Row 87
$ws.Range("H87").Value = 4750
$ws.Range("I87").Value = 0
$ws.Range("J87").Value = 4750
$ws.Range("K87").Value = 0
$ws.Range("L87").Value = 14250
$ws.Range("M87").ClearContents()
$ws.Range("N87").Value = -16746

# Row 90
$ws.Range("H90").Value = 4750
$ws.Range("I90").Value = 0
$ws.Range("J90").Value = 4750
$ws.Range("K90").Value = 0
$ws.Range("L90").Value = 42750
$ws.Range("M90").ClearContents()
$ws.Range("N90").Value = -55230

# Row 103
$ws.Range("H103").Value = 2239.4
$ws.Range("I103").Value = 299.5
$ws.Range("K103").Value = 898.5
$ws.Range("M103").Value = -19.5

# Row 107
$ws.Range("H107").Value = 517.58826
$ws.Range("J107").Value = 613.375
$ws.Range("L107").Value = 1840.125
$ws.Range("N107").Value = -5680.125

# Row 118
$ws.Range("H118").Value = 6349.75
$ws.Range("I118").Value = 6349.75
$ws.Range("K118").Value = 19049.25
$ws.Range("M118").Value = -17806.25

# Row 119
$ws.Range("H119").Value = 3443.3572
$ws.Range("I119").Value = 2361.9167
$ws.Range("K119").Value = 7085.750100000001
$ws.Range("M119").Value = -2247.750100000001

$ws = $wb.Worksheets.Item("GSM")
# Row 113
$ws.Range("H113").Value = 4049.625
$ws.Range("I113").Value = 3319
$ws.Range("J113").Value = 4780.25
$ws.Range("K113").Value = 3319
$ws.Range("L113").Value = 4780.25
$ws.Range("M113").Value = -1149
$ws.Range("N113").Value = -9120.25

# Row 122
$ws.Range("H122").Value = 2579.9429
$ws.Range("I122").Value = 2158.04
$ws.Range("K122").Value = 6474.12
$ws.Range("M122").Value = -4024.12

$ws = $wb.Worksheets.Item("LTW")
# Row 31
$ws.Range("H31").Value = 2967.8125
$ws.Range("J31").Value = 2452.5833
$ws.Range("L31").Value = 2452.5833
$ws.Range("N31").Value = -2948.5833

# Row 93
$ws.Range("H93").Value = 100012500
$ws.Range("I93").Value = 200000000
$ws.Range("J93").Value = 24994
$ws.Range("K93").Value = 200000000
$ws.Range("L93").Value = 24994
$ws.Range("M93").Value = -199998752
$ws.Range("N93").Value = -27490

# Row 100
$ws.Range("H100").Value = 76925624
$ws.Range("I100").Value = 500000800
$ws.Range("J100").Value = 2868.0908
$ws.Range("K100").Value = 500000800
$ws.Range("L100").Value = 2868.0908
$ws.Range("M100").Value = -500000259
$ws.Range("N100").Value = -3950.0908

$ws = $wb.Worksheets.Item("WVR")
# Row 107
$ws.Range("H107").Value = 599.5
$ws.Range("I107").Value = 599.5
$ws.Range("J107").Value = 0
$ws.Range("K107").Value = 1798.5
$ws.Range("L107").Value = 0
$ws.Range("M107").Value = 121.5
$ws.Range("N107").ClearContents()

# Row 108
$ws.Range("H108").Value = 0
$ws.Range("J108").Value = 0
$ws.Range("L108").Value = 0
$ws.Range("N108").ClearContents()

# Row 122
$ws.Range("H122").Value = 3803.7256
$ws.Range("I122").Value = 1895.2559
$ws.Range("J122").Value = 14061.75
$ws.Range("K122").Value = 5685.7677
$ws.Range("L122").Value = 42185.25
$ws.Range("M122").Value = -3235.7677
$ws.Range("N122").Value = -47085.25

# Row 123
$ws.Range("H123").Value = 93619.336
$ws.Range("J123").Value = 93619.336
$ws.Range("L123").Value = 93619.336
$ws.Range("N123").Value = -103419.336

